# timing analysis and plotting all subjects on same classification plot
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row 13: another "good dataset" trial block ---
$ws.Range("A13").Value = 20251023
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = 4
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 5
$ws.Range("F13").Value = 3
$ws.Range("G13").Value = 6

# --- New row 14: summary/classification counts for the new subject ---
# (string cells are written in the same order the original sharedStrings
#  table acquires them, so the new <si> entries land at the expected indices)
$ws.Range("B14").Value = 33
$ws.Range("F14").Value = "27,30,36,38"
$ws.Range("C14").Value = "19,24,32,28"
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = "10,16,17,18,20,35"
$ws.Range("G14").Value = "2,22,26,28,33,34"

# --- Column width adjustments (narrower, to fit all subjects on one plot/view) ---
$ws.Columns.Item(2).ColumnWidth = 19
$ws.Columns.Item(3).ColumnWidth = 17.1666666667
$ws.Columns.Item(4).ColumnWidth = 18.6666666667
$ws.Columns.Item(5).ColumnWidth = 19
$ws.Columns.Item(6).ColumnWidth = 20.6666666667
$ws.Columns.Item(7).ColumnWidth = 21.8333333333

# --- Move the active selection to reflect the new bottom of the data (below row 14) ---
$ws.Range("G15").Select()
